$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-06-18 Tuesday" "2024-06-19 Wednesday"
Replace-Text "87×42=" "50×33="
Replace-Text "30×19=" "71×79="
Replace-Text "26×30=" "20×16="
Replace-Text "88×69=" "50×82="
Replace-Text "64×59=" "83×41="
Replace-Text "75×53=" "17×71="
Replace-Text "29×68=" "50×37="
Replace-Text "75×66=" "64×91="
Replace-Text "18×59=" "50×96="
Replace-Text "85×87=" "92×31="
Replace-Text "52×85=" "13×20="
Replace-Text "19×75=" "97×65="
Replace-Text "98×74=" "27×98="
Replace-Text "78×61=" "58×93="
Replace-Text "84×84=" "53×56="
Replace-Text "70×26=" "69×11="
Replace-Text "20×25=" "79×58="
Replace-Text "80×35=" "34×68="
Replace-Text "50×16=" "61×49="
Replace-Text "86×99=" "18×27="
Replace-Text "15×18=" "53×24="
Replace-Text "64×54=" "35×57="
Replace-Text "20×34=" "88×63="
Replace-Text "12×62=" "16×23="
Replace-Text "89×73=" "91×33="
